$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Drop the old "Storico" (history) sheet entirely.
$wb.Worksheets("Storico").Delete() | Out-Null

# Rename the remaining "Flotta" sheet to "Sheet1".
$ws = $wb.Worksheets("Flotta")
$ws.Name = "Sheet1"

# Row 99 (targa GX834SK) is now marked as end-of-rental:
#   - operatore (column B) becomes "FINE RENT"
#   - a new column C cell records the change date as literal text "2025-12-18"
$ws.Range("B99").Value = "FINE RENT"

# Write "2025-12-18" into C99 as plain text (not an auto-converted date serial).
# A formula (string literal) evaluates to a text result; copying that result
# into C99 via paste-values keeps the destination a plain shared-string cell
# with no number-format/style attached, same as a freshly authored text cell.
$helper = $ws.Range("ZZ1")
$helper.Formula = "=""2025-12-18"""
$helper.Copy()
$ws.Range("C99").PasteSpecial(-4163)
$helper.Clear()
